$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.121.43"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "3.359.72"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'587.97"
$ws.Range("E5").Value = "  +6.34%  "
$ws.Range("D6").Value = "'188.45"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").Value = "'0.599"
$ws.Range("E7").Value = "  +3.56%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.185"
$ws.Range("E9").Value = "  +3.69%  "
$ws.Range("D10").Value = "'0.587"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").Value = "'47.58"
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("D12").Value = "'0.0000274"
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("D13").Value = "'658.06"
$ws.Range("E13").Value = "  +10.26%  "
$ws.Range("D14").Value = "3.900.16"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").Value = "'8.62"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "67.192.17"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.119"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'18.04"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "3.365.09"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "'11.22"
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("D21").Value = "'0.909"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").Value = "'18.07"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("D23").Value = "'5.11"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'100.82"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "'4.02"
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("D26").Value = "'2.84"
$ws.Range("E26").Value = "  +3.99%  "
$ws.Range("D27").Value = "'9.77"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").Value = "'32.20"
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("D29").Value = "'8.70"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "'6.91"
$ws.Range("E30").Value = "  +3.37%  "
$ws.Range("D31").Value = "'616.66"
$ws.Range("E31").Value = "  +8.44%  "
$ws.Range("D32").Value = "'3.95"
$ws.Range("E32").Value = "  +3.11%  "
$ws.Range("D33").Value = "'11.23"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("D35").Value = "3.882.80"
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'55.48"
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("D38").Value = "'2.82"
$ws.Range("E38").Value = "  +6.61%  "
$ws.Range("D39").Value = "'0.131"
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("D40").Value = "'33.73"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "'3.26"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").Value = "0.0₃0708"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("D43").Value = "'0.346"
$ws.Range("E43").Value = "  +3.10%  "
$ws.Range("D44").Value = "'3.39"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "'0.0423"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").Value = "'2.59"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").Value = "'2.86"
$ws.Range("E49").Value = "  -18.41%  "
$ws.Range("E50").Value = "  +8.81%  "
$ws.Range("D51").Value = "'129.75"
$ws.Range("E51").Value = "  +5.35%  "
